$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: "Jumlah" (E2) switches from a numeric value to a text value ---
# A leading apostrophe is how Excel marks a numeric-looking entry as text
# (quote-prefix) instead of letting it coerce to a number.
$ws.Range("E2").Value = "'60000"

# --- New payment-history rows 3-5 (riwayat_pembayaran) ---
# NIK, Plat, Nama, Tanggal_Bayar, Jumlah, Metode, Nama_Penerima, No_HP, Alamat, Jasa_Pengiriman
$rows = @(
    @("'1234456278949542", "BG4576HI", "Nia Rahmadani", "02-08-2025 16:24", "'65000", "Bank Mandiri", "Nia Rahmadani", "'085267947261", "Jl. Melati, Palembang", "JNE"),
    @("'1234456278949542", "BG4576HI", "Nia Rahmadani", "02-08-2025 16:31", "'65000", "Bank Mandiri", "Nia Rahmadani", "'085267947261", "Jl. Melati, Palembang", "JNE"),
    @("'1234456278949542", "BG4576HI", "Nia Rahmadani", "02-08-2025 16:45", 65000,    "Bank Negara Indonesia (BNI)", "Nia Rahmadani", "089012736819h", "Jl. Melati, Palembang", "J&T")
)

$startRow = 3
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $values = $rows[$i]
    for ($c = 1; $c -le $values.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $values[$c - 1]
    }
}
